# Commiting to Local.Changes for TestCase_F5
# Adds a new test case row (row 40) to the "Test Cases" sheet,
# mirroring the formatting of the existing rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy cell formatting down from the nearest matching rows so the new
# row's styles line up with the target (A/B/C/E match row 39's styles,
# D needs the plain bordered style used elsewhere, e.g. row 31).
$ws.Range("A39").Copy() | Out-Null
$ws.Range("A40").PasteSpecial(-4122) | Out-Null

$ws.Range("B39").Copy() | Out-Null
$ws.Range("B40").PasteSpecial(-4122) | Out-Null

$ws.Range("C39").Copy() | Out-Null
$ws.Range("C40").PasteSpecial(-4122) | Out-Null

$ws.Range("D31").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4122) | Out-Null

$ws.Range("E39").Copy() | Out-Null
$ws.Range("E40").PasteSpecial(-4122) | Out-Null

# Populate the new test case values.
$ws.Range("A40").Value = "AppreciateUnAppreciateOthersPost"
$ws.Range("B40").Value = "OPQA-342|OPQA-359"
$ws.Range("C40").Value = "Verify that user is able to Appreciate/Un Appreciate their others post"
$ws.Range("D40").Value = "Y"

# Move the active selection to the new row, as in the authored workbook.
$ws.Range("A40").Select() | Out-Null
